# Auto-generated COM-interop script applying the cryptos.xlsx price/volume update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.384.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.300.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.22%  "

$ws.Range("E7").Value = "  +0.66%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0911"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.50%  "

$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.970"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.647.88"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.296.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.502.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.23%  "

$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.10%  "

$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.05%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "277.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +20.08%  "

$ws.Range("E25").Value = "  -1.09%  "

$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("E32").Value = "  -2.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.137"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.117"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -11.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0368"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.94%  "

$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.14%  "

$ws.Range("B43").Value = "BitcoinSV"
$ws.Range("C43").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.78%  "

$ws.Range("E44").Value = "  -1.83%  "

$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "83.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.31%  "

$ws.Range("E47").Value = "  -0.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.76%  "

$ws.Range("E49").Value = "  -0.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.589.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.75%  "
